# ARKANSAS_2016.xlsx edit script
# - Removes the trailing metadata/footnote rows (964-968)
# - Renames the header row to short machine-friendly column names
# - Title-cases the lowercase Spanish connector words ("de", "del", "el",
#   "la", "los", "y") that appear inside the state/municipality names

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# 1) Drop the trailing metadata rows (sample size / source / author / date).
#    Row 963 never existed in the sheet, so the metadata starts at 964.
$ws.Range("A964:A968").EntireRow.Delete()

# 2) Rename the header row (row 1) to the new short column names.
$ws.Range("A1").Value = "mx_state"
$ws.Range("B1").Value = "mx_municipality"
$ws.Range("C1").Value = "n_matriculas"
$ws.Range("D1").Value = "pct_matriculas"

# 3) Title-case the Spanish connector words within the state (A) and
#    municipality (B) columns, for the data rows only (2-962).
$rng = $ws.Range("A2:B962")
$rng.Replace(" de ", " De ") | Out-Null
$rng.Replace(" del ", " Del ") | Out-Null
$rng.Replace(" el ", " El ") | Out-Null
$rng.Replace(" la ", " La ") | Out-Null
$rng.Replace(" los ", " Los ") | Out-Null
$rng.Replace(" y ", " Y ") | Out-Null
